$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checklist")

# Check the "Reset Zone" checkbox (linked to J13) and the checkbox linked to J20.
$ws.Range("J13").Value = $true
$ws.Range("J20").Value = $true

# Scroll the view so row 16 is the top-left visible row (topLeftCell="A16").
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Application.ActiveWindow.ScrollColumn = 1
